$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Helper: populate one "position statistics" sheet (rows 1-14, cols A-E)
# given the category label text (A), count (B), percentage (C) and the
# offset (D) for rows 1-13. Row 14 only carries A/B/C (D/E are blank on
# the template - F14 already holds the average formula from the copy).
# ------------------------------------------------------------------
function Set-PositionRow {
    param($ws, $row, $label, $b, $c, $d)
    $ws.Range("A$row").Value = $label
    $ws.Range("B$row").Value = $b
    $ws.Range("C$row").Value = $c
    if ($null -ne $d) {
        $ws.Range("D$row").Value = $d
    }
}

# ==================================================================
# 1. Add the three new sheets (20191015 / 20191017 / 20191021) by
#    copying the existing "20191014" sheet, which already carries the
#    right layout/styles/formulas (E column, F14 average, B15/C15
#    totals row).
# ==================================================================
$template = $wb.Worksheets.Item(2)

$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy($null, $last)
$sheet3 = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet3.Name = "20191015"

$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy($null, $last)
$sheet4 = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet4.Name = "20191017"

$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy($null, $last)
$sheet5 = $wb.Worksheets.Item($wb.Worksheets.Count)
$sheet5.Name = "20191021"

# ==================================================================
# 2. Fill in the per-day numbers for each new sheet.
# ==================================================================

# ---- 20191015 ----
Set-PositionRow $sheet3 1  "空仓"                               676  0.05 0
Set-PositionRow $sheet3 2  "0-9"                                302  0.02 5
Set-PositionRow $sheet3 3  "10-19"                               372  0.03 15
Set-PositionRow $sheet3 4  "20-29"                               457  0.03 25
Set-PositionRow $sheet3 5  "30-39"                               418  0.03 35
Set-PositionRow $sheet3 6  "40-49"                               525  0.04 45
Set-PositionRow $sheet3 7  "50-59"                               666  0.05 55
Set-PositionRow $sheet3 8  "60-69 (已选)"                        690  0.05 65
Set-PositionRow $sheet3 9  "70-79"                                732  0.06 75
Set-PositionRow $sheet3 10 "80-89"                                972  0.08 85
Set-PositionRow $sheet3 11 "90-99"                                1493 0.12 95
Set-PositionRow $sheet3 12 "100及以上（主动买入）"                 1443 0.12 100
Set-PositionRow $sheet3 13 "100及以上（被动买入，俗称被套死）"      727  0.06 100
Set-PositionRow $sheet3 14 "我是来给卫斯理打Call的"                2224 0.16 $null
$sheet3.Range("E1:E13").Formula = "=D1*B1"

# ---- 20191017 ----
Set-PositionRow $sheet4 1  "空仓"                               762  0.06 0
Set-PositionRow $sheet4 2  "0-9"                                356  0.02 5
Set-PositionRow $sheet4 3  "10-19"                               398  0.03 15
Set-PositionRow $sheet4 4  "20-29"                               477  0.03 25
Set-PositionRow $sheet4 5  "30-39"                               502  0.04 35
Set-PositionRow $sheet4 6  "40-49"                               560  0.04 45
Set-PositionRow $sheet4 7  "50-59"                               722  0.05 55
Set-PositionRow $sheet4 8  "60-69 (已选)"                        684  0.05 65
Set-PositionRow $sheet4 9  "70-79"                                813  0.06 75
Set-PositionRow $sheet4 10 "80-89"                                1013 0.08 85
Set-PositionRow $sheet4 11 "90-99"                                1531 0.12 95
Set-PositionRow $sheet4 12 "100及以上（主动买入）"                 1498 0.12 100
Set-PositionRow $sheet4 13 "100及以上（被动买入，俗称被套死）"      861  0.06 100
Set-PositionRow $sheet4 14 "我是来给卫斯理打Call的"                2206 0.17 $null
$sheet4.Range("E1:E13").Formula = "=D1*B1"

# ---- 20191021 ----
Set-PositionRow $sheet5 1  "空仓 (已选)"                        444  0.07 0
Set-PositionRow $sheet5 2  "0-9"                                173  0.02 5
Set-PositionRow $sheet5 3  "10-19"                               206  0.03 15
Set-PositionRow $sheet5 4  "20-29"                               221  0.03 25
Set-PositionRow $sheet5 5  "30-39"                               235  0.04 35
Set-PositionRow $sheet5 6  "40-49"                               240  0.04 45
Set-PositionRow $sheet5 7  "50-59"                               369  0.06 55
Set-PositionRow $sheet5 8  "60-69"                                294  0.05 65
Set-PositionRow $sheet5 9  "70-79"                                347  0.05 75
Set-PositionRow $sheet5 10 "80-89"                                432  0.07 85
Set-PositionRow $sheet5 11 "90-99"                                702  0.12 95
Set-PositionRow $sheet5 12 "100及以上（主动买入）"                 741  0.12 100
Set-PositionRow $sheet5 13 "100及以上（被动买入，俗称被套死）"      473  0.08 100
Set-PositionRow $sheet5 14 "我是来给卫斯理打Call的"                973  0.16 $null
$sheet5.Range("E1:E13").Formula = "=D1*B1"

# ---- 20191014 (existing sheet): re-assert E1:E13 as one shared
#      formula block, matching the newer sheets. ----
$template.Range("E1:E13").Formula = "=D1*B1"

# ==================================================================
# 3. Sheet-view / selection bookkeeping.
#    - 20191014 is no longer the active tab; its lingering C1:C14
#      selection moves to F14.
#    - 20191015 leaves the cursor at E20.
#    - 20191017 becomes the active tab, cursor at D24.
#    - 20191021 keeps an A1:C14 selection (mirrors 20191014's old one).
# ==================================================================
[void]$wb.Worksheets.Item("20191014").Range("F14").Select()
[void]$sheet3.Range("E20").Select()
[void]$sheet5.Range("A1:C14").Select()

$sheet4.Activate()
[void]$sheet4.Range("D24").Select()
